# Update the task description for Driver #100004 row: create basic CRUD
# functions for objects (previously described as Mop/extendedMop specific).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Add crud basic features for objects"
